# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
# (updated crypto price/volume figures scraped on Fri Apr 14 14:50:46 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.680.89'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.100.28'
$ws.Range("E3").Value = '  +5.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.94'
$ws.Range("E5").Value = '  +1.87%  '

$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("E7").Value = '  +3.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4326'
$ws.Range("E8").Value = '  +4.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08916'
$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.59'
$ws.Range("E10").Value = '  +8.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.166'
$ws.Range("E11").Value = '  +2.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.54'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.102.06'
$ws.Range("E13").Value = '  +5.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.698'
$ws.Range("E14").Value = '  +2.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.786'
$ws.Range("E15").Value = '  +4.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.79'
$ws.Range("E16").Value = '  +2.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.13%  '

$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("E19").Value = '  +2.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.91'
$ws.Range("E20").Value = '  -0.13%  '

$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.283'
$ws.Range("E22").Value = '  +1.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.722.49'
$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.30'
$ws.Range("E24").Value = '  +3.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.347.17'
$ws.Range("E25").Value = '  +5.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.283'
$ws.Range("E26").Value = '  +3.31%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.561'
$ws.Range("E28").Value = '  +5.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.94'
$ws.Range("E29").Value = '  -0.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.02'
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.203'
$ws.Range("E31").Value = '  +5.40%  '

$ws.Range("E32").Value = '  +2.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.137'
$ws.Range("E33").Value = '  +0.98%  '

$ws.Range("E34").Value = '  +16.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.837'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02591'
$ws.Range("E36").Value = '  +3.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.666'
$ws.Range("E37").Value = '  +7.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.516'
$ws.Range("E38").Value = '  +2.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06714'
$ws.Range("E39").Value = '  +1.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.59'

$ws.Range("E41").Value = '  +3.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6784'
$ws.Range("E42").Value = '  +2.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.252'
$ws.Range("E43").Value = '  +1.51%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.00'
$ws.Range("E45").Value = '  +1.91%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6376'
$ws.Range("E46").Value = '  +3.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.213'
$ws.Range("E47").Value = '  +0.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.629'
$ws.Range("E48").Value = '  -0.83%  '

$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.63'
$ws.Range("E50").Value = '  +2.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.196'
$ws.Range("E51").Value = '  +7.13%  '

